$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06191573587245055
$ws.Range("H2").Value = -4.379613200895871
$ws.Range("I2").Value = 459.6863173629243
$ws.Range("G3").Value = 0.06680946745287783
$ws.Range("H3").Value = -2.33258148145959
$ws.Range("G4").Value = -0.03908288381688925
$ws.Range("H4").Value = 12.94221537370796
$ws.Range("G5").Value = -0.02424925138234309
$ws.Range("H5").Value = 10.42629862114493
$ws.Range("G6").Value = -0.1025884088958581
$ws.Range("H6").Value = 3.256083002957712
$ws.Range("G7").Value = -0.08735638126406267
$ws.Range("H7").Value = 4.393931625934679
$ws.Range("G8").Value = -0.3653970310607456
$ws.Range("H8").Value = 0.4108693792994438
$ws.Range("G9").Value = -0.3754842595773724
$ws.Range("H9").Value = 3.744376429587384
$ws.Range("G10").Value = 0.02610195967060595
$ws.Range("H10").Value = 29.27523327306857
$ws.Range("G11").Value = 0.02448145428807504
$ws.Range("H11").Value = 7.874066972413285
$ws.Range("G12").Value = 0.2247337770452149
$ws.Range("H12").Value = 1.34667827098448
$ws.Range("G13").Value = 0.2153840120456268
$ws.Range("H13").Value = -4.361675256710567
$ws.Range("G14").Value = -0.04615002856872588
$ws.Range("H14").Value = -9.60673551989068
$ws.Range("G15").Value = -0.05210792298475388
$ws.Range("H15").Value = -9.249411017350521
$ws.Range("G16").Value = 0.2212101149799407
$ws.Range("H16").Value = 4.06293367722375
$ws.Range("G17").Value = 0.2263566474584515
$ws.Range("H17").Value = 2.633750473735192
$ws.Range("G18").Value = 0.08061226410872555
$ws.Range("H18").Value = 10.39714941470925
$ws.Range("G19").Value = 0.07571973855551552
$ws.Range("H19").Value = 0.5066841462666574
$ws.Range("G20").Value = -0.08136331454456529
$ws.Range("H20").Value = -8.510496476315566
$ws.Range("G21").Value = -0.07781129014068847
$ws.Range("H21").Value = 10.11627522863095
$ws.Range("G22").Value = 0.06630947889539876
$ws.Range("H22").Value = -9.785350104596994
$ws.Range("G23").Value = 0.07163371307421811
$ws.Range("H23").Value = 4.832918341066388
$ws.Range("G24").Value = 0.06506498021889888
$ws.Range("H24").Value = -2.322626440380759
$ws.Range("G25").Value = 0.05937974298847255
$ws.Range("H25").Value = 8.396649347406695
$ws.Range("G26").Value = 0.1174438028484849
$ws.Range("H26").Value = -1.596352567984715
$ws.Range("G27").Value = 0.1196586471448159
$ws.Range("H27").Value = 5.092617702514871
$ws.Range("G28").Value = 0.1309651687587763
$ws.Range("H28").Value = 1.324367587146582
$ws.Range("G29").Value = 0.1446308546695255
$ws.Range("H29").Value = -4.116641104503471
$ws.Range("G30").Value = 0.09058289824644754
$ws.Range("H30").Value = 7.443709865561764
$ws.Range("G31").Value = 0.08964818302980701
$ws.Range("H31").Value = 9.744557270502037
$ws.Range("G32").Value = 0.06384434676086118
$ws.Range("H32").Value = 19.64855675782632
$ws.Range("G33").Value = 0.063815768477714
$ws.Range("H33").Value = 15.51847560769001
$ws.Range("G34").Value = 0.02009991780780899
$ws.Range("H34").Value = 15.80132004517278
$ws.Range("G35").Value = 0.02415252286088796
$ws.Range("H35").Value = 42.90851748760415
$ws.Range("G36").Value = -0.01986601364123035
$ws.Range("H36").Value = 31.60428807373076
$ws.Range("G37").Value = -0.0195917605294944
$ws.Range("H37").Value = 41.1019961879947
$ws.Range("G38").Value = 0.07444800409273465
$ws.Range("H38").Value = -4.898550732797786
$ws.Range("G39").Value = 0.07524739617794246
$ws.Range("H39").Value = -3.212721774277396
$ws.Range("G40").Value = 0.06003681220024818
$ws.Range("H40").Value = -9.322322462120972
$ws.Range("G41").Value = 0.06074698932327535
$ws.Range("H41").Value = -6.575920350032629
$ws.Range("G42").Value = 0.08924260629315311
$ws.Range("H42").Value = 14.7205589533018
$ws.Range("G43").Value = 0.09240221731531657
$ws.Range("H43").Value = 15.26576786867596
$ws.Range("G44").Value = 0.08554526179043372
$ws.Range("H44").Value = -3.061895648305034
$ws.Range("G45").Value = 0.08786365760164518
$ws.Range("H45").Value = -2.787951183339356
$ws.Range("G46").Value = -0.003501737187738383
$ws.Range("H46").Value = -27.97711280117426
$ws.Range("G47").Value = 0.0121446370673871
$ws.Range("H47").Value = 25521.33909979539
$ws.Range("G48").Value = -0.09606006067442592
$ws.Range("H48").Value = 0.04691068970927208
$ws.Range("G49").Value = -0.09929202972256911
$ws.Range("H49").Value = 9.385407677800906
$ws.Range("G50").Value = 0.1656810221171058
$ws.Range("H50").Value = -2.825975874365572
$ws.Range("G51").Value = 0.1755880022628731
$ws.Range("H51").Value = 3.387486731159944
$ws.Range("G52").Value = 0.06205750454069868
$ws.Range("H52").Value = -12.54811664314424
$ws.Range("G53").Value = 0.0725440818558547
$ws.Range("H53").Value = 12.79302184868881
$ws.Range("G54").Value = -0.1281564478033565
$ws.Range("H54").Value = -0.2766699938493567
$ws.Range("G55").Value = -0.1293170370406249
$ws.Range("H55").Value = -11.02922699069733
$ws.Range("G56").Value = 0.1956237438727999
$ws.Range("H56").Value = 2.948576004348846
$ws.Range("G57").Value = 0.195965168288365
$ws.Range("H57").Value = -1.479629785011756
